$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Data" sheet: insert 3 new rows (2024, 2023, 2022) at the top of the
#    time series (right after the header row), pushing the existing years
#    (2021..2006) down by three rows.
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Rows.Item(2).Insert()
$wsData.Rows.Item(2).Insert()
$wsData.Rows.Item(2).Insert()

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $wsData.Range("A2") "2024"
$wsData.Range("B2").Value = 0.95

Set-TextValue $wsData.Range("A3") "2023"
$wsData.Range("B3").Value = 0.95

Set-TextValue $wsData.Range("A4") "2022"
$wsData.Range("B4").Value = 0.88

# ---------------------------------------------------------------------------
# 2) "Metadata" sheet: insert a new "actualizacion" / "Julio 2025" row
#    between "observaciones" and "cita", and refresh the "observaciones"
#    and "cita" texts.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$nuevaObservacion = "Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. En julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH. Durante el año 2020 y hasta julio del año 2021 se suspende el relevamiento de la información necesaria para construir indicadores relativos al nivel y la trayectoria educativa. A partir de esta fecha, las preguntas se relevan en el formulario presencial. Un conjunto importante de indicadores educativos tienen un efecto estacional, por lo que no se recomienda comparar los resultados del segundo semestre del 2021 con la información anual. Las estimaciones desde 2022 se calculan a partir de la muestra de implantación."

$wsMeta.Range("B8").Value = $nuevaObservacion

# Insert the new "actualizacion" row right after "observaciones" (row 8),
# so it becomes row 9 and everything below shifts down by one row.
$wsMeta.Rows.Item(9).Insert()
$wsMeta.Range("A9").Value = "actualizacion"
$wsMeta.Range("B9").Value = "Julio 2025"

# "cita" text (row 10 now) gains a trailing newline.
$nuevaCita = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE`n"
$wsMeta.Range("B10").Value = $nuevaCita
